# Reorder the two worksheets so "review_info" comes before "hotel_info"
# (matches the sheet-order swap seen in xl/workbook.xml), then add the
# new "State" column (value "Louisiana") to the "hotel_info" sheet,
# right after the "Hotel_Name" column.

$wb = $excel.ActiveWorkbook

$hotel  = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# Move "hotel_info" to sit right after "review_info" -> order becomes
# review_info, hotel_info.
$hotel.Move($null, $review)

# Re-fetch (defensive - names didn't change, just tab order) and insert
# the new "State" column (C) in hotel_info, shifting City/Zip/etc. right.
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"
